$wb = $excel.ActiveWorkbook

# --- Rename the second sheet (encapsulation of the race-enrollment data model) ---
$wsGeneral = $wb.Worksheets.Item("General_Enrollment")
$wsRace = $wb.Worksheets.Item("Enrollment_by_Race")
$wsRace.Name = "Race_Enrollment"

# --- Fix the bug with selecting the wrong sparse matrix / wrong active sheet ---
# Previously sheet2 ("Enrollment_by_Race") was the active/selected tab with
# selection on F10. Make sheet1 ("General_Enrollment") the active tab instead,
# with a selection on O14, and clear the selection left on the other sheet.

$wsRace.Range("M24").Select()

$wsGeneral.Activate()
$wsGeneral.Range("O14").Select()

# --- New internal data model / column layout on General_Enrollment ---
# Add explicit column widths for the new non-freshman / non-first-year /
# non-first-time columns (K, M, O).
$wsGeneral.Columns.Item(11).ColumnWidth = 11.67
$wsGeneral.Columns.Item(13).ColumnWidth = 11.67
$wsGeneral.Columns.Item(15).ColumnWidth = 14.2
